$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename benchmarker "NREL" -> "NLR" for all data rows that used it (rows 2-22)
$ws.Range("A2:A22").Value = "NLR"

# Update the example_benchmarker rows (23-25): MPI Tasks (E) and the
# Epsilon Total/Benchmark Time columns (H/J) were corrected.
$ws.Range("E23").Value = 2
$ws.Range("H23").Value = 300
$ws.Range("J23").Value = 305

$ws.Range("E24").Value = 4
$ws.Range("H24").Value = 200
$ws.Range("J24").Value = 205

$ws.Range("E25").Value = 8
$ws.Range("H25").Value = 100
$ws.Range("J25").Value = 105

# Update the active selection to match the saved view state
$ws.Range("G22").Select()
